# Update the roster table on the active sheet to match the new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data rows (Player, Position, Team) for rows 2..18.
$data = @(
    @("Jalen Green",        "PG,SG",    "Houston Rockets"),
    @("Russell Westbrook",  "PG,SG",    "Denver Nuggets"),
    @("Payton Pritchard",   "PG,SG",    "Boston Celtics"),
    @("Jalen Suggs",        "PG,SG",    "Orlando Magic"),
    @("Ayo Dosunmu",        "PG,SG,SF", "Chicago Bulls"),
    @("Paolo Banchero",     "SF,PF",    "Orlando Magic"),
    @("Jaylen Brown",       "SG,SF",    "Boston Celtics"),
    @("Kyle Kuzma",         "PF",       "Washington Wizards"),
    @("Nikola Jokic",       "C",        "Denver Nuggets"),
    @("Rudy Gobert",        "C",        "Minnesota Timberwolves"),
    @("Dejounte Murray",    "PG,SG",    "New Orleans Pelicans"),
    @("Jakob Poeltl",       "C",        "Toronto Raptors"),
    @("Chris Paul",         "PG",       "San Antonio Spurs"),
    @("Deni Avdija",        "SF,PF",    "Portland Trail Blazers"),
    @("Pascal Siakam",      "SF,PF,C",  "Indiana Pacers"),
    @("Chet Holmgren",      "PF,C",     "Oklahoma City Thunder"),
    @("Khris Middleton",    "SF",       "Milwaukee Bucks")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# The new table only has 17 data rows (rows 2-18); the old sheet had 18
# rows (2-19), so delete the now-unused last row entirely.
$ws.Rows.Item(19).Delete()
